$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5782.0557
$ws.Range("I12").Value = 124
$ws.Range("J12").Value = 14673.286
$ws.Range("K12").Value = 124
$ws.Range("L12").Value = 14673.286
$ws.Range("M12").Value = 46
$ws.Range("N12").Value = -15013.286
$ws.Range("H33").Value = 688.16
$ws.Range("I33").Value = 558.5789
$ws.Range("K33").Value = 558.5789
$ws.Range("M33").Value = -329.5789
$ws.Range("H55").Value = 947.4
$ws.Range("I55").Value = 1158.9
$ws.Range("J55").Value = 524.4
$ws.Range("K55").Value = 1158.9
$ws.Range("L55").Value = 524.4
$ws.Range("M55").Value = -944.9000000000001
$ws.Range("N55").Value = -952.4
$ws.Range("H62").Value = 2812.2144
$ws.Range("I62").Value = 2277.5
$ws.Range("K62").Value = 2277.5
$ws.Range("M62").Value = -1653.5
$ws.Range("H65").Value = 2812.2144
$ws.Range("I65").Value = 2277.5
$ws.Range("K65").Value = 11387.5
$ws.Range("M65").Value = -8267.5
$ws.Range("H92").Value = 448
$ws.Range("I92").Value = 416.5263
$ws.Range("J92").Value = 597.5
$ws.Range("K92").Value = 416.5263
$ws.Range("L92").Value = 597.5
$ws.Range("M92").Value = 831.4737
$ws.Range("N92").Value = -3093.5
$ws.Range("H98").Value = 2117.4866
$ws.Range("J98").Value = 5099.5
$ws.Range("L98").Value = 5099.5
$ws.Range("N98").Value = -8095.5
$ws.Range("H100").Value = 3835.1052
$ws.Range("I100").Value = 3749.0833
$ws.Range("J100").Value = 3982.5715
$ws.Range("K100").Value = 3749.0833
$ws.Range("L100").Value = 3982.5715
$ws.Range("M100").Value = -3208.0833
$ws.Range("N100").Value = -5064.5715
$ws.Range("H106").Value = 2284.2104
$ws.Range("I106").Value = 1488.8889
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1488.8889
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -857.8888999999999
$ws.Range("N106").Value = -4262
$ws.Range("H115").Value = 993.2143
$ws.Range("I115").Value = 418.63635
$ws.Range("K115").Value = 1255.90905
$ws.Range("M115").Value = 311.09095
$ws.Range("H122").Value = 2117.4866
$ws.Range("J122").Value = 5099.5
$ws.Range("L122").Value = 15298.5
$ws.Range("N122").Value = -20198.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 77025.8
$ws.Range("I53").Value = 25000
$ws.Range("J53").Value = 90032.25
$ws.Range("K53").Value = 25000
$ws.Range("L53").Value = 90032.25
$ws.Range("M53").Value = -24318
$ws.Range("N53").Value = -91396.25
$ws.Range("H61").Value = 3493.7058
$ws.Range("I61").Value = 2429.1538
$ws.Range("J61").Value = 6953.5
$ws.Range("K61").Value = 2429.1538
$ws.Range("L61").Value = 6953.5
$ws.Range("M61").Value = -2217.1538
$ws.Range("N61").Value = -7377.5
$ws.Range("H113").Value = 35000
$ws.Range("J113").Value = 35000
$ws.Range("L113").Value = 35000
$ws.Range("N113").Value = -43678
$ws.Range("H136").Value = 3493.7058
$ws.Range("I136").Value = 2429.1538
$ws.Range("J136").Value = 6953.5
$ws.Range("K136").Value = 7287.4614
$ws.Range("L136").Value = 20860.5
$ws.Range("M136").Value = -4737.4614
$ws.Range("N136").Value = -25960.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 28509
$ws.Range("J27").Value = 28509
$ws.Range("L27").Value = 28509
$ws.Range("N27").Value = -28893
$ws.Range("H82").Value = 17690.268
$ws.Range("J82").Value = 29206
$ws.Range("L82").Value = 29206
$ws.Range("N82").Value = -29972
$ws.Range("H85").Value = 17690.268
$ws.Range("J85").Value = 29206
$ws.Range("L85").Value = 29206
$ws.Range("N85").Value = -31858
$ws.Range("H94").Value = 838.11536
$ws.Range("I94").Value = 585.4286
$ws.Range("J94").Value = 1899.4
$ws.Range("K94").Value = 585.4286
$ws.Range("L94").Value = 1899.4
$ws.Range("M94").Value = -134.4286
$ws.Range("N94").Value = -2801.4
$ws.Range("H97").Value = 13483.333
$ws.Range("I97").Value = 6966.6665
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 6966.6665
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -5975.6665
$ws.Range("N97").Value = -21982
$ws.Range("H107").Value = 2378.5
$ws.Range("I107").Value = 2166.5
$ws.Range("J107").Value = 2537.5
$ws.Range("K107").Value = 2166.5
$ws.Range("L107").Value = 2537.5
$ws.Range("M107").Value = -246.5
$ws.Range("N107").Value = -6377.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1636.1818
$ws.Range("I16").Value = 746.6667
$ws.Range("K16").Value = 746.6667
$ws.Range("M16").Value = -459.6667
$ws.Range("H28").Value = 35643
$ws.Range("J28").Value = 35643
$ws.Range("L28").Value = 35643
$ws.Range("N28").Value = -36133
$ws.Range("H31").Value = 3191.1409
$ws.Range("I31").Value = 2172.077
$ws.Range("J31").Value = 4433.125
$ws.Range("K31").Value = 2172.077
$ws.Range("L31").Value = 4433.125
$ws.Range("M31").Value = -1877.077
$ws.Range("N31").Value = -5023.125
$ws.Range("H34").Value = 3191.1409
$ws.Range("I34").Value = 2172.077
$ws.Range("J34").Value = 4433.125
$ws.Range("K34").Value = 2172.077
$ws.Range("L34").Value = 4433.125
$ws.Range("M34").Value = -1970.077
$ws.Range("N34").Value = -4837.125
$ws.Range("H43").Value = 32000
$ws.Range("J43").Value = 32000
$ws.Range("L43").Value = 32000
$ws.Range("N43").Value = -32368
$ws.Range("H88").Value = 31114.334
$ws.Range("J88").Value = 31114.334
$ws.Range("L88").Value = 31114.334
$ws.Range("N88").Value = -31926.334
$ws.Range("H91").Value = 31114.334
$ws.Range("J91").Value = 31114.334
$ws.Range("L91").Value = 31114.334
$ws.Range("N91").Value = -33922.334
$ws.Range("H95").Value = 25608
$ws.Range("J95").Value = 25608
$ws.Range("L95").Value = 25608
$ws.Range("N95").Value = -31100
$ws.Range("H101").Value = 32000
$ws.Range("J101").Value = 32000
$ws.Range("L101").Value = 32000
$ws.Range("N101").Value = -38490
$ws.Range("H113").Value = 1636.1818
$ws.Range("I113").Value = 746.6667
$ws.Range("K113").Value = 746.6667
$ws.Range("M113").Value = 1423.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1372.3077
$ws.Range("I68").Value = 473.33334
$ws.Range("J68").Value = 2142.8572
$ws.Range("K68").Value = 1420.00002
$ws.Range("L68").Value = 6428.571599999999
$ws.Range("M68").Value = -609.0000199999999
$ws.Range("N68").Value = -8050.571599999999
$ws.Range("H71").Value = 1372.3077
$ws.Range("I71").Value = 473.33334
$ws.Range("J71").Value = 2142.8572
$ws.Range("K71").Value = 4260.00006
$ws.Range("L71").Value = 19285.7148
$ws.Range("M71").Value = -204.0000600000003
$ws.Range("N71").Value = -27397.7148
$ws.Range("H82").Value = 2400
$ws.Range("H85").Value = 2400
$ws.Range("H86").Value = 1600
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = 586
$ws.Range("N86").Value = -11372
$ws.Range("H89").Value = 1600
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 27000
$ws.Range("M89").Value = 4128
$ws.Range("N89").Value = -38856
$ws.Range("H97").Value = 912.63635
$ws.Range("J97").Value = 1529.8
$ws.Range("L97").Value = 4589.4
$ws.Range("N97").Value = -5581.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4399.6587
$ws.Range("I132").Value = 4841.1665
$ws.Range("K132").Value = 14523.4995
$ws.Range("M132").Value = -11993.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1001.1579
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1161.4667
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1161.4667
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1751.4667
$ws.Range("H27").Value = 1001.1579
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1161.4667
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1161.4667
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1375.4667
$ws.Range("H46").Value = 5875
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 7666.6665
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 7666.6665
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -8042.6665
$ws.Range("H68").Value = 4001.5715
$ws.Range("I68").Value = 1037
$ws.Range("J68").Value = 6225
$ws.Range("K68").Value = 1037
$ws.Range("L68").Value = 6225
$ws.Range("M68").Value = -288
$ws.Range("N68").Value = -7723
$ws.Range("H71").Value = 4001.5715
$ws.Range("I71").Value = 1037
$ws.Range("J71").Value = 6225
$ws.Range("K71").Value = 5185
$ws.Range("L71").Value = 31125
$ws.Range("M71").Value = -1441
$ws.Range("N71").Value = -38613
$ws.Range("H122").Value = 3171.4644
$ws.Range("I122").Value = 2542.25
$ws.Range("J122").Value = 3643.375
$ws.Range("K122").Value = 7626.75
$ws.Range("L122").Value = 10930.125
$ws.Range("M122").Value = -5176.75
$ws.Range("N122").Value = -15830.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1937.12
$ws.Range("J122").Value = 3859.6
$ws.Range("L122").Value = 11578.8
$ws.Range("N122").Value = -16478.8
